$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 7.572506307224296
$ws.Range("G2").Value = 3415
$ws.Range("H2").Value = 5.765304948399171
$ws.Range("I2").Value = 1.313461538461538
$ws.Range("J2").Value = 193.0788461538461
$ws.Range("K2").Value = 1113.158427161972
$ws.Range("L2").Value = 847.4998274146782
$ws.Range("M2").Value = 1.25627879636851
$ws.Range("F3").Value = 7.55268244677915
$ws.Range("G3").Value = 3409
$ws.Range("H3").Value = 5.760332754950364
$ws.Range("I3").Value = 1.311153846153846
$ws.Range("J3").Value = 165.2053846153846
$ws.Range("K3").Value = 951.6379882941729
$ws.Range("L3").Value = 725.8019271237459
$ws.Range("M3").Value = 1.36552498637767
$ws.Range("F4").Value = 7.50043214450328
$ws.Range("G4").Value = 3414
$ws.Range("H4").Value = 5.712104152228625
$ws.Range("I4").Value = 1.313076923076923
$ws.Range("J4").Value = 1.313076923076923
$ws.Range("K4").Value = 7.50043214450328
$ws.Range("L4").Value = 5.712104152228625
$ws.Range("M4").Value = 1.87585807934027
$ws.Range("F5").Value = 7.534238085261641
$ws.Range("G5").Value = 3414
$ws.Range("H5").Value = 5.737849742730013
$ws.Range("I5").Value = 1.313076923076923
$ws.Range("J5").Value = 1.313076923076923
$ws.Range("K5").Value = 7.534238085261641
$ws.Range("L5").Value = 5.737849742730013
$ws.Range("M5").Value = 1.98376488784939
$ws.Range("F6").Value = 8.878141583495536
$ws.Range("G6").Value = 3417
$ws.Range("H6").Value = 6.755390142548549
$ws.Range("I6").Value = 1.314230769230769
$ws.Range("J6").Value = 173.4784615384615
$ws.Range("K6").Value = 1171.914689021411
$ws.Range("L6").Value = 891.7114988164085
$ws.Range("M6").Value = 1.837775307783576
$ws.Range("F7").Value = 9.024545432902086
$ws.Range("G7").Value = 3482
$ws.Range("H7").Value = 6.738603712103798
$ws.Range("I7").Value = 1.339230769230769
$ws.Range("J7").Value = 196.8669230769231
$ws.Range("K7").Value = 1326.608178636607
$ws.Range("L7").Value = 990.5747456792583
$ws.Range("M7").Value = 1.497172087318456
$ws.Range("F8").Value = 9.006205375083971
$ws.Range("G8").Value = 3482
$ws.Range("H8").Value = 6.724909240441793
$ws.Range("I8").Value = 1.339230769230769
$ws.Range("J8").Value = 1.339230769230769
$ws.Range("K8").Value = 9.006205375083971
$ws.Range("L8").Value = 6.724909240441793
$ws.Range("M8").Value = 2.75319698316317
$ws.Range("F9").Value = 8.850115772952966
$ws.Range("G9").Value = 3409
$ws.Range("H9").Value = 6.749868292659934
$ws.Range("I9").Value = 1.311153846153846
$ws.Range("J9").Value = 103.5811538461538
$ws.Range("K9").Value = 699.1591460632843
$ws.Range("L9").Value = 533.2395951201348
$ws.Range("M9").Value = 1.375307991116891
$ws.Range("F10").Value = 11.16145205226164
$ws.Range("G10").Value = 3491
$ws.Range("H10").Value = 8.312739998819898
$ws.Range("I10").Value = 1.342692307692308
$ws.Range("J10").Value = 1.342692307692308
$ws.Range("K10").Value = 11.16145205226164
$ws.Range("L10").Value = 8.312739998819898
$ws.Range("M10").Value = 3.159807075995271
$ws.Range("F11").Value = 10.89486911554118
$ws.Range("G11").Value = 3413
$ws.Range("H11").Value = 8.299636595489909
$ws.Range("I11").Value = 1.312692307692308
$ws.Range("J11").Value = 99.76461538461538
$ws.Range("K11").Value = 828.0100527811294
$ws.Range("L11").Value = 630.7723812572331
$ws.Range("M11").Value = 1.726836754813277
$ws.Range("F12").Value = 12.21462677222745
$ws.Range("G12").Value = 3812
$ws.Range("H12").Value = 8.331067578119454
$ws.Range("I12").Value = 1.466153846153846
$ws.Range("J12").Value = 115.8261538461538
$ws.Range("K12").Value = 964.9555150059682
$ws.Range("L12").Value = 658.1543386714368
$ws.Range("M12").Value = 1.898153000404145
$ws.Range("F13").Value = 12.21315217615776
$ws.Range("G13").Value = 3817
$ws.Range("H13").Value = 8.319150028297139
$ws.Range("I13").Value = 1.468076923076923
$ws.Range("J13").Value = 111.5738461538461
$ws.Range("K13").Value = 928.1995653879899
$ws.Range("L13").Value = 632.2554021505825
$ws.Range("M13").Value = 1.935784619921005
$ws.Range("F14").Value = 13.08832600787552
$ws.Range("G14").Value = 4122
$ws.Range("H14").Value = 8.255615628451324
$ws.Range("I14").Value = 1.585384615384615
$ws.Range("J14").Value = 1.585384615384615
$ws.Range("K14").Value = 13.08832600787552
$ws.Range("L14").Value = 8.255615628451324
$ws.Range("M14").Value = 3.273390334569668
$ws.Range("F15").Value = 11.44023702710331
$ws.Range("G15").Value = 4268
$ws.Range("H15").Value = 6.969216558216638
$ws.Range("I15").Value = 1.641538461538462
$ws.Range("J15").Value = 129.6815384615385
$ws.Range("K15").Value = 903.7787251411615
$ws.Range("L15").Value = 550.5681080991144
$ws.Range("M15").Value = 1.777812834011855
$ws.Range("F16").Value = 11.54427496778571
$ws.Range("G16").Value = 4329
$ws.Range("H16").Value = 6.933498479150579
$ws.Range("I16").Value = 1.665
$ws.Range("J16").Value = 1.665
$ws.Range("K16").Value = 11.54427496778571
$ws.Range("L16").Value = 6.933498479150579
$ws.Range("M16").Value = 2.226890641285864
$ws.Range("F17").Value = 11.41800602768974
$ws.Range("G17").Value = 4267
$ws.Range("H17").Value = 6.957303883757518
$ws.Range("I17").Value = 1.641153846153846
$ws.Range("J17").Value = 124.7276923076923
$ws.Range("K17").Value = 867.7684581044205
$ws.Range("L17").Value = 528.7550951655714
$ws.Range("M17").Value = 1.809753955388824
$ws.Range("F18").Value = 11.38806946575303
$ws.Range("G18").Value = 4260
$ws.Range("H18").Value = 6.950464932149737
$ws.Range("I18").Value = 1.638461538461538
$ws.Range("J18").Value = 1.638461538461538
$ws.Range("K18").Value = 11.38806946575303
$ws.Range("L18").Value = 6.950464932149737
$ws.Range("M18").Value = 2.848156173384833
$ws.Range("F19").Value = 12.60236766102066
$ws.Range("G19").Value = 4641
$ws.Range("H19").Value = 7.060149950151633
$ws.Range("I19").Value = 1.785
$ws.Range("J19").Value = 1.785
$ws.Range("K19").Value = 12.60236766102066
$ws.Range("L19").Value = 7.060149950151633
$ws.Range("M19").Value = 3.318203405146741
$ws.Range("F20").Value = 13.61299186526415
$ws.Range("G20").Value = 4242
$ws.Range("H20").Value = 8.34365366564988
$ws.Range("I20").Value = 1.631538461538462
$ws.Range("J20").Value = 128.8915384615385
$ws.Range("K20").Value = 1075.426357355868
$ws.Range("L20").Value = 659.1486395863406
$ws.Range("M20").Value = 2.115458935862049
$ws.Range("F21").Value = 13.74589299891787
$ws.Range("G21").Value = 4319
$ws.Range("H21").Value = 8.274906644405293
$ws.Range("I21").Value = 1.661153846153846
$ws.Range("J21").Value = 1.661153846153846
$ws.Range("K21").Value = 13.74589299891787
$ws.Range("L21").Value = 8.274906644405293
$ws.Range("M21").Value = 2.651582759491257
$ws.Range("F22").Value = 13.65883553855037
$ws.Range("G22").Value = 4243
$ws.Range("H22").Value = 8.369779024329709
$ws.Range("I22").Value = 1.631923076923077
$ws.Range("J22").Value = 1.631923076923077
$ws.Range("K22").Value = 13.65883553855037
$ws.Range("L22").Value = 8.369779024329709
$ws.Range("M22").Value = 3.866816340963609
$ws.Range("F23").Value = 13.57196208277016
$ws.Range("G23").Value = 4236
$ws.Range("H23").Value = 8.33028834164363
$ws.Range("I23").Value = 1.629230769230769
$ws.Range("J23").Value = 123.8215384615385
$ws.Range("K23").Value = 1031.469118290532
$ws.Range("L23").Value = 633.1019139649159
$ws.Range("M23").Value = 2.15115599011907
$ws.Range("F24").Value = 14.63407654354711
$ws.Range("G24").Value = 4592
$ws.Range("H24").Value = 8.285844732844618
$ws.Range("I24").Value = 1.766153846153846
$ws.Range("J24").Value = 1.766153846153846
$ws.Range("K24").Value = 14.63407654354711
$ws.Range("L24").Value = 8.285844732844618
$ws.Range("M24").Value = 4.473637199362352
$ws.Range("J25").Value = 1.423109760062423
$ws.Range("K25").Value = 10.13351681494935
$ws.Range("L25").Value = 7.120685346508236
$ws.Range("M25").Value = 53.17831914003776
$ws.Range("N25").Value = 2600
$ws.Range("P25").Value = 0.02020916923218286

$ws.Range("Q25").Value = "(51.07192742410905, 55.28471085596654)"
$ws.Range("R25").Value = "(50.40991859910287, 55.94671968097271)"
